$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Update the Start Year value from 2016 to 2020 (next year fuel price fix)
$ws.Range("B1").Value = 2020

# Update the selected/active cell on the sheet to D5 to match the saved view state
$ws.Activate()
$ws.Range("D5").Select()
